$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old hyperlinks (and their relationships) before rewriting the table
$ws.Hyperlinks.Delete()

# Clear out the old data rows (2-10); new rows will be written below
$ws.Range("A2:F10").ClearContents()

$data = @(
    ,@("infrastructure", $null, 121359, "https://www.wikidata.org/wiki/Q121359", "infrastructure", 3)
    ,@("invoice", $null, 190581, "https://www.wikidata.org/wiki/Q190581", "bill", 2)
    ,@("receivership", "ORG", 474341, "https://www.wikidata.org/wiki/Q474341", "Administration", 1)
    ,@("economy", $null, 159810, "https://www.wikidata.org/wiki/Q159810", "economy", 1)
    ,@("United States of America", "GPE", 30, "https://www.wikidata.org/wiki/Q30", "America", 1)
    ,@("repair", $null, 2144962, "https://www.wikidata.org/wiki/Q2144962", "repairs", 1)
    ,@("bridge", $null, 12280, "https://www.wikidata.org/wiki/Q12280", "bridges", 1)
    ,@("Michigan", "GPE", 1166, "https://www.wikidata.org/wiki/Q1166", "Michigan", 1)
    ,@("road", $null, 34442, "https://www.wikidata.org/wiki/Q34442", "roads", 1)
    ,@("percent", "PERCENT", 11229, "https://www.wikidata.org/wiki/Q11229", "%", 1)
    ,@("pork", $null, 191768, "https://www.wikidata.org/wiki/Q191768", "pork", 1)
    ,@("Democratic Party", "NORP", 29552, "https://www.wikidata.org/wiki/Q29552", "Democrats", 1)
    ,@("productivity", $null, 2111958, "https://www.wikidata.org/wiki/Q2111958", "productivity", 1)
    ,@("hundred", "CARDINAL", 313354, "https://www.wikidata.org/wiki/Q313354", "hundreds", 1)
    ,@("global warming", $null, 7942, "https://www.wikidata.org/wiki/Q7942", "climate change", 1)
    ,@("American Capitalism", $null, 4743290, "https://www.wikidata.org/wiki/Q4743290", "American capitalism", 1)
    ,@("agenda", $null, 1758159, "https://www.wikidata.org/wiki/Q1758159", "agenda", 1)
    ,@("First Step", $null, 3269975, "https://www.wikidata.org/wiki/Q3269975", "first step", 1)
    ,@("Green New Deal", $null, 2068307, "https://www.wikidata.org/wiki/Q2068307", "Green New Deal", 1)
    ,@("Americans", "ORG", 846570, "https://www.wikidata.org/wiki/Q846570", "Americans", 1)
    ,@("Houston Texans", "ORG", 223514, "https://www.wikidata.org/wiki/Q223514", "Texans", 1)
    ,@("act", $null, 820655, "https://www.wikidata.org/wiki/Q820655", "act", 1)
    ,@("Immediately", $null, 6004788, "https://www.wikidata.org/wiki/Q6004788", "IMMEDIATELY", 1)
    ,@("2012 VP113", "ORG", 15980607, "https://www.wikidata.org/wiki/Q15980607", "Biden", 1)
    ,@("Texas", "GPE", 1439, "https://www.wikidata.org/wiki/Q1439", "Texas", 1)
    ,@("illegal immigration", $null, 856681, "https://www.wikidata.org/wiki/Q856681", "illegal immigrants", 1)
    ,@("generation", $null, 213381, "https://www.wikidata.org/wiki/Q213381", "generations", 1)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne $null) {
        $ws.Cells.Item($r, 2).Value = $row[1]
    }
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 4), $row[3])
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Re-apply the Hyperlink cell style after all links are created (avoids
# the engine allocating a fresh style xf per Hyperlinks.Add call).
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 4).Style = "Hyperlink"
}
